# 04 - prog_taller.xlsx — "Se modifican algunas fechas para vestir mejor
# los datos." (a handful of workshop-session dates in column A were
# retyped by hand; every date after the touched cells re-derives through
# the existing "=A(row-1)+7" weekly-increment chain, and the SQL-insert
# text built in column H recalculates from those dates automatically.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("usuario")

# Directly retyped dates (these replace the previous formula/value in
# place, exactly like a user typing a new date into the cell). Written
# as the underlying 1900-date-system serial numbers so the stored cell
# is a plain literal, matching "2023-03-14" etc.
$ws.Range("A4").Value  = 44999   # 2023-03-14
$ws.Range("A5").Value  = 45008   # 2023-03-23
$ws.Range("A8").Value  = 45030   # 2023-04-14
$ws.Range("A9").Value  = 45033   # 2023-04-17
$ws.Range("A11").Value = 45048   # 2023-05-02
$ws.Range("A12").Value = 45056   # 2023-05-10
$ws.Range("A13").Value = 45066   # 2023-05-20
$ws.Range("A18").Value = 45100   # 2023-06-23

# Selection / scroll position left where the editor ended up.
$ws.Range("A19").Select() | Out-Null
